$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.735.19'
$ws.Range("E2").Value = '  +4.43%  '
$ws.Range("D3").Value = '3.344.33'
$ws.Range("E3").Value = '  +4.37%  '
$ws.Range("D5").Value = "'561.64"
$ws.Range("E5").Value = '  +4.46%  '
$ws.Range("D6").Value = "'152.52"
$ws.Range("E6").Value = '  +4.57%  '
$ws.Range("D7").Value = "'0.999"
$ws.Range("E8").Value = '  +0.18%  '
$ws.Range("D9").Value = "'7.44"
$ws.Range("E9").Value = '  +0.96%  '
$ws.Range("E10").Value = '  +3.79%  '
$ws.Range("E11").Value = '  +0.04%  '
$ws.Range("D12").Value = '3.921.98'
$ws.Range("E12").Value = '  +4.31%  '
$ws.Range("E13").Value = '  +0.23%  '
$ws.Range("D14").Value = "'26.91"
$ws.Range("E14").Value = '  +2.77%  '
$ws.Range("E15").Value = '  +3.05%  '
$ws.Range("D16").Value = '62.738.18'
$ws.Range("E16").Value = '  +4.25%  '
$ws.Range("D17").Value = '3.316.94'
$ws.Range("E17").Value = '  +3.20%  '
$ws.Range("D18").Value = "'6.34"
$ws.Range("E18").Value = '  +1.22%  '
$ws.Range("D19").Value = "'13.83"
$ws.Range("E19").Value = '  +4.43%  '
$ws.Range("E20").Value = '  +0.30%  '
$ws.Range("D21").Value = "'384.00"
$ws.Range("E21").Value = '  +1.12%  '
$ws.Range("D22").Value = "'0.998"
$ws.Range("E22").Value = '  -0.22%  '
$ws.Range("D23").Value = "'0.535"
$ws.Range("E23").Value = '  +1.82%  '
$ws.Range("D24").Value = "'70.14"
$ws.Range("E24").Value = '  -0.02%  '
$ws.Range("E25").Value = '  +5.15%  '
$ws.Range("D26").Value = "'8.92"
$ws.Range("E26").Value = '  -0.36%  '
$ws.Range("D27").Value = '0.0₃0948'
$ws.Range("E27").Value = '  +4.81%  '
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = '  +0.17%  '
$ws.Range("D29").Value = "'6.57"
$ws.Range("E29").Value = '  +5.54%  '
$ws.Range("E30").Value = '  +3.76%  '
$ws.Range("B31").Value = 'EthereumClassic'
$ws.Range("C31").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D31").Value = "'22.89"
$ws.Range("E31").Value = '  +2.44%  '
$ws.Range("B32").Value = 'NEARProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D32").Value = "'5.57"
$ws.Range("E32").Value = '  +2.43%  '
$ws.Range("E33").Value = '  +7.13%  '
$ws.Range("E34").Value = '  +0.41%  '
$ws.Range("E35").Value = '  +1.91%  '
$ws.Range("E36").Value = '  +7.74%  '
$ws.Range("E37").Value = '  +12.46%  '
$ws.Range("D38").Value = "'26.98"
$ws.Range("E38").Value = '  +5.22%  '
$ws.Range("D39").Value = "'0.0739"
$ws.Range("E39").Value = '  +4.72%  '
$ws.Range("D40").Value = '2.807.28'
$ws.Range("E40").Value = '  +0.30%  '
$ws.Range("E41").Value = '  +6.18%  '
$ws.Range("B42").Value = 'OKB'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D42").Value = "'40.55"
$ws.Range("E42").Value = '  +1.72%  '
$ws.Range("B43").Value = 'Mantle'
$ws.Range("C43").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D43").Value = "'0.744"
$ws.Range("E43").Value = '  +3.64%  '
$ws.Range("D44").Value = "'4.25"
$ws.Range("E44").Value = '  -0.12%  '
$ws.Range("E45").Value = '  +3.59%  '
$ws.Range("D46").Value = '3.389.11'
$ws.Range("E46").Value = '  +4.30%  '
$ws.Range("D47").Value = "'21.94"
$ws.Range("E47").Value = '  +6.00%  '
$ws.Range("E48").Value = '  -1.67%  '
$ws.Range("E49").Value = '  +1.83%  '
$ws.Range("D50").Value = "'288.00"
$ws.Range("E50").Value = '  +6.00%  '
$ws.Range("E51").Value = '  -1.31%  '
